$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 26147
$ws.Range("E2").Value = 508944419962
$ws.Range("F2").Value = 8354720331
$ws.Range("G2").Value = 0.94137

$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 1672.63
$ws.Range("E3").Value = 201082244499
$ws.Range("F3").Value = 8654979325
$ws.Range("G3").Value = 0.79074

$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 0.999888
$ws.Range("E4").Value = 82782523125
$ws.Range("F4").Value = 12887084993
$ws.Range("G4").Value = -0.07586

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 216.55
$ws.Range("E5").Value = 33314509412
$ws.Range("F5").Value = 352761642
$ws.Range("G5").Value = 0.63412

$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "XRP"
$ws.Range("D6").Value = 0.5439580000000001
$ws.Range("E6").Value = 28719636842
$ws.Range("F6").Value = 1339261159
$ws.Range("G6").Value = 8.025550000000001

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "USD Coin"
$ws.Range("D7").Value = 0.99897
$ws.Range("E7").Value = 25963465049
$ws.Range("F7").Value = 4338512583
$ws.Range("G7").Value = -0.11487

$ws.Range("B8").Value = "STETH"
$ws.Range("C8").Value = "Lido Staked Ether"
$ws.Range("D8").Value = 1673.48
$ws.Range("E8").Value = 13737122850
$ws.Range("F8").Value = 4876107
$ws.Range("G8").Value = 0.84137

$ws.Range("B9").Value = "ADA"
$ws.Range("C9").Value = "Cardano"
$ws.Range("D9").Value = 0.267765
$ws.Range("E9").Value = 9378891535
$ws.Range("F9").Value = 145815909
$ws.Range("G9").Value = 1.47134

$ws.Range("B10").Value = "DOGE"
$ws.Range("C10").Value = "Dogecoin"
$ws.Range("D10").Value = 0.064161
$ws.Range("E10").Value = 9017665408
$ws.Range("F10").Value = 248625097
$ws.Range("G10").Value = 1.2268

$ws.Range("B11").Value = "SOL"
$ws.Range("C11").Value = "Solana"
$ws.Range("D11").Value = 21.89
$ws.Range("E11").Value = 8928837932
$ws.Range("F11").Value = 226571899
$ws.Range("G11").Value = 1.06256

$ws.Range("B12").Value = "TRX"
$ws.Range("C12").Value = "TRON"
$ws.Range("D12").Value = 0.074889
$ws.Range("E12").Value = 6696107467
$ws.Range("F12").Value = 155580362
$ws.Range("G12").Value = 1.51578

$ws.Range("B13").Value = "DOT"
$ws.Range("C13").Value = "Polkadot"
$ws.Range("D13").Value = 4.52
$ws.Range("E13").Value = 5725168908
$ws.Range("F13").Value = 81404260
$ws.Range("G13").Value = 0.96538

$ws.Range("B14").Value = "MATIC"
$ws.Range("C14").Value = "Polygon"
$ws.Range("D14").Value = 0.576191
$ws.Range("E14").Value = 5368123016
$ws.Range("F14").Value = 172996725
$ws.Range("G14").Value = -0.41416

$ws.Range("B15").Value = "SHIB"
$ws.Range("C15").Value = "Shiba Inu"
$ws.Range("D15").Value = 0.000008389999999999999
$ws.Range("E15").Value = 4945250637
$ws.Range("F15").Value = 150797169
$ws.Range("G15").Value = -1.08194

$ws.Range("B16").Value = "TON"
$ws.Range("C16").Value = "Toncoin"
$ws.Range("D16").Value = 1.4
$ws.Range("E16").Value = 4802898610
$ws.Range("F16").Value = 15630116
$ws.Range("G16").Value = 4.96089

$ws.Range("B17").Value = "LTC"
$ws.Range("C17").Value = "Litecoin"
$ws.Range("D17").Value = 64.83
$ws.Range("E17").Value = 4773195935
$ws.Range("F17").Value = 420143451
$ws.Range("G17").Value = 1.57336

$ws.Range("B18").Value = "WBTC"
$ws.Range("C18").Value = "Wrapped Bitcoin"
$ws.Range("D18").Value = 26226
$ws.Range("E18").Value = 4256375499
$ws.Range("F18").Value = 58600219
$ws.Range("G18").Value = 1.07334

$ws.Range("B19").Value = "DAI"
$ws.Range("C19").Value = "Dai"
$ws.Range("D19").Value = 0.999447
$ws.Range("E19").Value = 3933070634
$ws.Range("F19").Value = 59512314
$ws.Range("G19").Value = -0.01357

$ws.Range("B20").Value = "AVAX"
$ws.Range("C20").Value = "Avalanche"
$ws.Range("D20").Value = 10.86
$ws.Range("E20").Value = 3736125562
$ws.Range("F20").Value = 79838820
$ws.Range("G20").Value = 1.35132

$ws.Range("B21").Value = "UNI"
$ws.Range("C21").Value = "Uniswap"
$ws.Range("D21").Value = 4.9
$ws.Range("E21").Value = 3695707344
$ws.Range("F21").Value = 81427803
$ws.Range("G21").Value = -0.00773

$ws.Range("B22").Value = "BCH"
$ws.Range("C22").Value = "Bitcoin Cash"
$ws.Range("D22").Value = 188.22
$ws.Range("E22").Value = 3665294031
$ws.Range("F22").Value = 165327543
$ws.Range("G22").Value = -0.12432

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "LEO Token"
$ws.Range("D23").Value = 3.85
$ws.Range("E23").Value = 3614542293
$ws.Range("F23").Value = 492404
$ws.Range("G23").Value = -1.96984

$ws.Range("B24").Value = "XLM"
$ws.Range("C24").Value = "Stellar"
$ws.Range("D24").Value = 0.130233
$ws.Range("E24").Value = 3570084986
$ws.Range("F24").Value = 134743513
$ws.Range("G24").Value = 9.23034

$ws.Range("B25").Value = "LINK"
$ws.Range("C25").Value = "Chainlink"
$ws.Range("D25").Value = 6.2
$ws.Range("E25").Value = 3335883953
$ws.Range("F25").Value = 140922148
$ws.Range("G25").Value = 0.80148

$ws.Range("B26").Value = "BUSD"
$ws.Range("C26").Value = "Binance USD"
$ws.Range("D26").Value = 0.998811
$ws.Range("E26").Value = 3240223616
$ws.Range("F26").Value = 498448592
$ws.Range("G26").Value = -0.17773

$ws.Range("B27").Value = "TUSD"
$ws.Range("C27").Value = "TrueUSD"
$ws.Range("D27").Value = 0.999126
$ws.Range("E27").Value = 2761696957
$ws.Range("F27").Value = 1161988410
$ws.Range("G27").Value = 0.07542

$ws.Range("B28").Value = "XMR"
$ws.Range("C28").Value = "Monero"
$ws.Range("D28").Value = 148.19
$ws.Range("E28").Value = 2684835985
$ws.Range("F28").Value = 50372267
$ws.Range("G28").Value = 3.04446

$ws.Range("B29").Value = "OKB"
$ws.Range("C29").Value = "OKB"
$ws.Range("D29").Value = 43.83
$ws.Range("E29").Value = 2626366118
$ws.Range("F29").Value = 4992821
$ws.Range("G29").Value = 3.49681

$ws.Range("B30").Value = "ATOM"
$ws.Range("C30").Value = "Cosmos Hub"
$ws.Range("D30").Value = 7.84
$ws.Range("E30").Value = 2292591070
$ws.Range("F30").Value = 93820497
$ws.Range("G30").Value = 3.63317

$ws.Range("B31").Value = "ETC"
$ws.Range("C31").Value = "Ethereum Classic"
$ws.Range("D31").Value = 15.54
$ws.Range("E31").Value = 2215647967
$ws.Range("F31").Value = 55613456
$ws.Range("G31").Value = 0.87935

$ws.Range("B32").Value = "HBAR"
$ws.Range("C32").Value = "Hedera"
$ws.Range("D32").Value = 0.06232
$ws.Range("E32").Value = 2056081744
$ws.Range("F32").Value = 90202973
$ws.Range("G32").Value = -4.27712

$ws.Range("B33").Value = "FIL"
$ws.Range("C33").Value = "Filecoin"
$ws.Range("D33").Value = 3.56
$ws.Range("E33").Value = 1574386555
$ws.Range("F33").Value = 63378654
$ws.Range("G33").Value = 0.57764

$ws.Range("B34").Value = "ICP"
$ws.Range("C34").Value = "Internet Computer"
$ws.Range("D34").Value = 3.56
$ws.Range("E34").Value = 1571248282
$ws.Range("F34").Value = 13634734
$ws.Range("G34").Value = 1.66302

$ws.Range("B35").Value = "LDO"
$ws.Range("C35").Value = "Lido DAO"
$ws.Range("D35").Value = 1.66
$ws.Range("E35").Value = 1458687561
$ws.Range("F35").Value = 30772007
$ws.Range("G35").Value = 0.64093

$ws.Range("B36").Value = "QNT"
$ws.Range("C36").Value = "Quant"
$ws.Range("D36").Value = 100.01
$ws.Range("E36").Value = 1454294861
$ws.Range("F36").Value = 12262212
$ws.Range("G36").Value = -0.6747

$ws.Range("B37").Value = "MNT"
$ws.Range("C37").Value = "Mantle"
$ws.Range("D37").Value = 0.426022
$ws.Range("E37").Value = 1377793336
$ws.Range("F37").Value = 9373096
$ws.Range("G37").Value = -0.15404

$ws.Range("B38").Value = "CRO"
$ws.Range("C38").Value = "Cronos"
$ws.Range("D38").Value = 0.052297
$ws.Range("E38").Value = 1371820858
$ws.Range("F38").Value = 2899323
$ws.Range("G38").Value = 0.64296

$ws.Range("B39").Value = "APT"
$ws.Range("C39").Value = "Aptos"
$ws.Range("D39").Value = 6.04
$ws.Range("E39").Value = 1364650353
$ws.Range("F39").Value = 71113112
$ws.Range("G39").Value = 1.10033

$ws.Range("B40").Value = "ARB"
$ws.Range("C40").Value = "Arbitrum"
$ws.Range("D40").Value = 1.036
$ws.Range("E40").Value = 1318623451
$ws.Range("F40").Value = 129140637
$ws.Range("G40").Value = 2.4892

$ws.Range("B41").Value = "VET"
$ws.Range("C41").Value = "VeChain"
$ws.Range("D41").Value = 0.01625167
$ws.Range("E41").Value = 1180954868
$ws.Range("F41").Value = 25085614
$ws.Range("G41").Value = 2.93322

$ws.Range("B42").Value = "NEAR"
$ws.Range("C42").Value = "NEAR Protocol"
$ws.Range("D42").Value = 1.17
$ws.Range("E42").Value = 1100934464
$ws.Range("F42").Value = 43223163
$ws.Range("G42").Value = 1.64873

$ws.Range("B43").Value = "OP"
$ws.Range("C43").Value = "Optimism"
$ws.Range("D43").Value = 1.47
$ws.Range("E43").Value = 1053828862
$ws.Range("F43").Value = 66391038
$ws.Range("G43").Value = 0.13124

$ws.Range("B44").Value = "MKR"
$ws.Range("C44").Value = "Maker"
$ws.Range("D44").Value = 1110.31
$ws.Range("E44").Value = 999837767
$ws.Range("F44").Value = 42632921
$ws.Range("G44").Value = 1.89615

$ws.Range("B45").Value = "RETH"
$ws.Range("C45").Value = "Rocket Pool ETH"
$ws.Range("D45").Value = 1823.58
$ws.Range("E45").Value = 913435101
$ws.Range("F45").Value = 10297184
$ws.Range("G45").Value = 0.97265

$ws.Range("B46").Value = "GRT"
$ws.Range("C46").Value = "The Graph"
$ws.Range("D46").Value = 0.09526900000000001
$ws.Range("E46").Value = 870904701
$ws.Range("F46").Value = 51708501
$ws.Range("G46").Value = 3.07448

$ws.Range("B47").Value = "KAS"
$ws.Range("C47").Value = "Kaspa"
$ws.Range("D47").Value = 0.04144989
$ws.Range("E47").Value = 841880643
$ws.Range("F47").Value = 10174393
$ws.Range("G47").Value = 2.67073

$ws.Range("B48").Value = "AAVE"
$ws.Range("C48").Value = "Aave"
$ws.Range("D48").Value = 56.94
$ws.Range("E48").Value = 826532196
$ws.Range("F48").Value = 57731717
$ws.Range("G48").Value = 1.68619

$ws.Range("B49").Value = "XDC"
$ws.Range("C49").Value = "XDC Network"
$ws.Range("D49").Value = 0.058346
$ws.Range("E49").Value = 808337205
$ws.Range("F49").Value = 9862155
$ws.Range("G49").Value = 2.50754

$ws.Range("B50").Value = "FRAX"
$ws.Range("C50").Value = "Frax"
$ws.Range("D50").Value = 0.997945
$ws.Range("E50").Value = 805507599
$ws.Range("F50").Value = 7701598
$ws.Range("G50").Value = 0.02373

$ws.Range("B51").Value = "WBT"
$ws.Range("C51").Value = "WhiteBIT Coin"
$ws.Range("D51").Value = 5.42
$ws.Range("E51").Value = 779302478
$ws.Range("F51").Value = 11509550
$ws.Range("G51").Value = 0.08634
